# Update the cached "datetimeFigureOut" footer date field text from
# 5/12/25 -> 5/23/25 everywhere it is used: the slide master, every
# custom (slide) layout, and the notes master.

$p = $ppt.ActivePresentation
$oldDate = "5/12/25"
$newDate = "5/23/25"

function Update-DateShapes($container) {
    $shapes = $container.Shapes
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.HasTextFrame) {
            if ($sh.TextFrame.HasText) {
                $tr = $sh.TextFrame.TextRange
                if ($tr.Text -eq $oldDate) {
                    $tr.Text = $newDate
                }
            }
        }
    }
}

# Slide master footer date placeholder.
Update-DateShapes($p.SlideMaster)

# Every custom layout inherits / overrides its own footer date placeholder.
$layouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    Update-DateShapes($layouts.Item($li))
}

# Notes master footer date placeholder.
Update-DateShapes($p.NotesMaster)
